$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value2 = 111234523
$ws.Range("B2").Value2 = 90666
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value2 = 4364
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "Dropptaggsvamp"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "Hydnellum ferrugineum"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q2").Value2 = 374912.3124199872
$ws.Range("R2").Value2 = 6871174.188302284
$ws.Range("S2").Value2 = 5
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = "18:52"
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "18:52"

# --- Row 4 ---
$ws.Range("A4").Value2 = 111232569
$ws.Range("B4").Value2 = 78081
$ws.Range("E4").Value2 = 229821
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "Vedflamlav"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Ramboldia elabens"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = ""
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = ""
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "Smolbäcken , Dlr"
$ws.Range("Q4").Value2 = 374784.9101014594
$ws.Range("R4").Value2 = 6871121.416580504
$ws.Range("S4").Value2 = 4
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "17:35"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "17:35"
$ws.Range("AC4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()

# --- Row 5 ---
$ws.Range("A5").Value2 = 111233768
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "Smolbäcken, Dlr"
$ws.Range("Q5").Value2 = 374893.0173954847
$ws.Range("R5").Value2 = 6871124.034136346
$ws.Range("S5").Value2 = 4
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "18:20"
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "18:20"

# --- Row 6 ---
$ws.Range("A6").Value2 = 111233604
$ws.Range("B6").Value2 = 56398
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value2 = 100109
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1"
$ws.Range("L6").NumberFormat = "@"
$ws.Range("L6").Value = ""
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "äldre spår"
$ws.Range("Q6").Value2 = 374894.0328106415
$ws.Range("R6").Value2 = 6871113.164681672
$ws.Range("S6").Value2 = 5
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = "18:15"
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "18:15"
$ws.Range("AC6").NumberFormat = "@"
$ws.Range("AC6").Value = "Hackmärken på tall"
$ws.Range("J6").ClearContents()

# --- Row 10 ---
$ws.Range("A10").Value2 = 111348142
$ws.Range("B10").Value2 = 89646
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value2 = 65
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "Fläckporing"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "Anthoporia albobrunnea"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "(Romell) Karasiński & Niemelä"
$ws.Range("Q10").Value2 = 374924.2261946606
$ws.Range("R10").Value2 = 6871177.518043431
$ws.Range("Z10").NumberFormat = "@"
$ws.Range("Z10").Value = "18:00"
$ws.Range("AB10").NumberFormat = "@"
$ws.Range("AB10").Value = "18:00"
$ws.Range("AC10").NumberFormat = "@"
$ws.Range("AC10").Value = "Växer på undersidan av kolad stubbe."
$ws.Range("AI10").NumberFormat = "@"
$ws.Range("AI10").Value = "Tallnaturskog"

# --- Row 11 ---
$ws.Range("A11").Value2 = 111346661
$ws.Range("B11").Value2 = 90854
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value2 = 2079
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "Nordtagging"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "Odonticium romellii"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "(S.Lundell) Parmasto"
$ws.Range("Q11").Value2 = 374943.5121711227
$ws.Range("R11").Value2 = 6871150.435782712
$ws.Range("Z11").NumberFormat = "@"
$ws.Range("Z11").Value = "19:00"
$ws.Range("AB11").NumberFormat = "@"
$ws.Range("AB11").Value = "19:00"
$ws.Range("AC11").NumberFormat = "@"
$ws.Range("AC11").Value = "Växer under kolad tallåga."
$ws.Range("AI11").NumberFormat = "@"
$ws.Range("AI11").Value = "Lavhedtallskog"
